# Update the cryptocurrency price/volume table on Sheet1 (rows 2-51, columns B-E)
# with the latest snapshot of coin name, link, price and 1h volume change.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row, Coin, Link, Price, Volume(1h)
$rows = @(
    @(2, "Bitcoin", "https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc", "31.575.36", "  +5.80%  "),
    @(3, "Ethereum", "https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth", "1.710.51", "  +4.19%  "),
    @(4, "TetherUSD", "https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt", "0.999", "  -0.09%  "),
    @(5, "BNB", "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb", "222.34", "  +3.12%  "),
    @(6, "XRP", "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp", "0.536", "  +3.09%  "),
    @(7, "USDC", "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc", "0.999", "  -0.06%  "),
    @(8, "Solana", "https://coinranking.com/coin/zNZHO_Sjf+solana-sol", "29.90", "  +3.42%  "),
    @(9, "Cardano", "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada", "0.270", "  +3.46%  "),
    @(10, "Dogecoin", "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge", "0.0648", "  +6.32%  "),
    @(11, "TRON", "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx", "0.0911", "  +1.15%  "),
    @(12, "WrappedliquidstakedEther2.0", "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth", "1.955.73", "  +4.28%  "),
    @(13, "WrappedEther", "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth", "1.715.63", "  +4.70%  "),
    @(14, "Polygon", "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic", "0.613", "  +3.64%  "),
    @(15, "Chainlink", "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link", "10.24", "  +7.88%  "),
    @(16, "Polkadot", "https://coinranking.com/coin/25W7FG7om+polkadot-dot", "4.21", "  +8.29%  "),
    @(17, "WrappedBTC", "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc", "31.555.29", "  +5.67%  "),
    @(18, "Litecoin", "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc", "67.22", "  +4.16%  "),
    @(19, "BitcoinCash", "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch", "251.52", "  +4.63%  "),
    @(20, "ShibaInu", "https://coinranking.com/coin/xz24e0BjL+shibainu-shib", "0.0₃0726", "  +3.11%  "),
    @(21, "Dai", "https://coinranking.com/coin/MoTuySvg7+dai-dai", "0.999", "  -0.13%  "),
    @(22, "Avalanche", "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax", "10.16", "  +2.15%  "),
    @(23, "Uniswap", "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni", "4.25", "  +2.57%  "),
    @(24, "Toncoin", "https://coinranking.com/coin/67YlI0K1b+toncoin-ton", "2.16", "  -1.39%  "),
    @(25, "Monero", "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr", "159.49", "  +1.31%  "),
    @(26, "EthereumClassic", "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc", "16.06", "  +3.19%  "),
    @(27, "Stellar", "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm", "0.113", "  +2.93%  "),
    @(28, "Cosmos", "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom", "6.80", "  +2.30%  "),
    @(29, "BinanceUSD", "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd", "0.999", "  -0.12%  "),
    @(30, "Filecoin", "https://coinranking.com/coin/ymQub4fuB+filecoin-fil", "3.83", "  +13.00%  "),
    @(31, "Hedera", "https://coinranking.com/coin/jad286TjB+hedera-hbar", "0.0504", "  +1.74%  "),
    @(32, "PancakeSwap", "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake", "1.15", "  +3.82%  "),
    @(33, "InternetComputer(DFINITY)", "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp", "3.38", "  +5.55%  "),
    @(34, "Maker", "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr", "1.526.27", "  +7.22%  "),
    @(35, "LidoDAOToken", "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo", "1.74", "  +2.25%  "),
    @(36, "TrustWalletToken", "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt", "1.04", "  +2.17%  "),
    @(37, "Aave", "https://coinranking.com/coin/ixgUfzmLR+aave-aave", "83.05", "  +8.47%  "),
    @(38, "ImmutableX", "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx", "0.609", "  +7.95%  "),
    @(39, "VeChain", "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet", "0.0181", "  +4.73%  "),
    @(40, "MXToken", "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx", "2.71", "  +0.18%  "),
    @(41, "HuobiToken", "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht", "2.31", "  +0.70%  "),
    @(42, "ARBITRUM", "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb", "0.855", "  +2.45%  "),
    @(43, "RenderToken", "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr", "2.04", "  +4.69%  "),
    @(44, "Kaspa", "https://coinranking.com/coin/V8GxkwWow+kaspa-kas", "0.0505", "  +0.74%  "),
    @(45, "WEMIXToken", "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix", "1.03", "  +3.19%  "),
    @(46, "PaxDollar", "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp", "1.00", "  -0.03%  "),
    @(47, "BitcoinSV", "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv", "52.25", "  +6.17%  "),
    @(48, "FraxShare", "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs", "5.59", "  +4.18%  "),
    @(49, "RocketPoolETH", "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth", "1.847.27", "  +3.59%  "),
    @(50, "BabyDogeCoin", "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge", "0.0₆0120", "  +10.42%  "),
    @(51, "Quant", "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt", "93.75", "  +0.35%  ")
)

foreach ($entry in $rows) {
    $r = $entry[0]
    $coin = $entry[1]
    $link = $entry[2]
    $price = $entry[3]
    $volume = $entry[4]

    $ws.Range("B$r").Value = $coin
    $ws.Range("C$r").Value = $link

    # The Price column holds plain text (e.g. "31.575.36" or "0.999"). Some of these
    # values would otherwise be auto-recognized as numbers by Excel, so force text by
    # pre-pending an apostrophe, then reset the style back to Normal so no stray
    # quote-prefix / number-format style is left behind on the cell.
    if ($price -match '^-?\d+(\.\d+)?$') {
        $ws.Range("D$r").Value = "'" + $price
        $ws.Range("D$r").Style = "Normal"
    } else {
        $ws.Range("D$r").Value = $price
    }

    $ws.Range("E$r").Value = $volume
}
